# Update the cryptocurrency price/volume data per the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.670.37'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.700.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '677.18'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.73'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.147'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.09'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.66'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.709.15'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '69.634.38'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '16.07'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.49'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '471.80'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.81'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.01%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '80.47'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.846.18'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000126'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.94'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.11'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.53%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.00%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.96'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.689.59'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.163'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.48'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.21'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.69%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.23'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0902'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '166.98'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.943'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '47.01'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.76'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.36'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000278'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'SuiNetwork'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.11'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.30'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.88'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.44%  '
